$d = $word.ActiveDocument

$oldText = "多云转小雨，今天学习了多分支管理，创建了一个dev分支。"
$newText = "使用git创建分支简单又快捷。"

# Locate the paragraph that currently holds the diary entry text.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*$oldText*") {
        $targetIndex = $i
        break
    }
}

# Insert a brand-new paragraph right before the target one, seeded from
# the paragraph that precedes it (this mirrors the formatting - rFonts
# hint="eastAsia" in the paragraph mark - used by every other entry in
# the diary) and carry the original sentence over into it.
$prevPara = $d.Paragraphs.Item($targetIndex - 1)
$prevPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($targetIndex)
$newPara.Range.Text = $oldText

# The paragraph that used to hold the original sentence (now shifted down
# by one) gets the new "simple and facile" sentence instead.
$target = $d.Paragraphs.Item($targetIndex + 1)
$target.Range.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                            $true, 1, $false, $newText, 2)

Write-Output "done"
